# Apply the LOM3228.xlsx content restructuring described by the commit diff.
# The worksheet holds a label/value/value table (columns A/B/C) describing a
# course. The block of rows 13-23 had its label/value pairs reshuffled (some
# removed, some values relocated to different rows), and the trailing row 24
# was removed entirely, shrinking the sheet from A1:C24 down to A1:C23.
#
# Notes on technique:
#  - Column A always uses style index 1 (bold), column B style index 2
#    (wrap), column C style index 3 (red wrap). When a brand new cell is
#    created in column B of a row that already has a column-A cell, it can
#    otherwise incorrectly inherit column A's style, so newly-created B/C
#    cells have their formatting explicitly re-applied by copying format
#    (not value) from an existing, correctly-styled B/C cell (B3/C3).
#  - "01/01/2012" looks like a date, so assigning it directly would get
#    auto-converted into a date serial number with date formatting. It is
#    routed through a scratch cell (Z1) pre-formatted as Text so the literal
#    string is preserved, then normal wrap formatting is restored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "Programa resumido:" / "Semestral" / "Semestral" -------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: "Short syllabus:" only ---------------------------------------
$ws.Range("B14:C14").Clear()
$ws.Range("A14").Value = "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: "Programa:" / "01/01/2012" / "01/01/2012" --------------------
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "01/01/2012"
$ws.Range("Z1").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("B3").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Programa:"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: "Syllabus:" only ----------------------------------------------
$ws.Range("B16:C16").Clear()
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: "Avaliação:" only, default height -----------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# --- Row 18: "Método:" / "6495737 - Durval Rodrigues Junior" x2 ------------
$ws.Range("B3").Copy($ws.Range("B18"))
$ws.Range("C3").Copy($ws.Range("C18"))
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C18").Value = "6495737 - Durval Rodrigues Junior"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: "Critério:" / lab-experiments description ---------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."
$ws.Range("C19").Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: "Norma de recuperação:" / arithmetic-mean grading description -
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3"
$ws.Range("C20").Value = "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3"
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: "Bibliografia:" / make-up exam description ---------------------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: "Requisitos:" only, default height -----------------------------
$ws.Range("B22:C22").Clear()
$ws.Range("A22").Value = "Requisitos:"
$ws.Rows.Item(22).AutoFit()

# --- Row 23: requisite course text (cols B/C only; A23 cleared) ------------
$ws.Range("A23").Clear()
$ws.Range("B3").Copy($ws.Range("B23"))
$ws.Range("C3").Copy($ws.Range("C23"))
$ws.Range("B23").Value = "LOB1019 -  Física II  (Requisito)`n"
$ws.Range("C23").Value = "LOB1019 -  Física II  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- Remove the now-obsolete row 24 entirely (shrinks dimension to C23) ---
$ws.Rows.Item(24).Delete()
